$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the first paragraph whose text (paragraph mark included)
# equals the given target string, starting the search at paragraph index
# $startIndex (1-based). Returns the Paragraph COM object.
# ---------------------------------------------------------------------------
function Find-ParagraphByText($doc, $targetText, $startIndex) {
    $count = $doc.Paragraphs.Count
    for ($i = $startIndex; $i -le $count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text -eq $targetText) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Insert nine new bullet paragraphs right after the very first paragraph
#    ("Encoding:") and before "URN: Encoding. Semantic Hashing."
# ---------------------------------------------------------------------------
$anchor1 = Find-ParagraphByText $d "Encoding:`r" 1
if ($anchor1 -eq $null) {
    throw "Could not locate anchor paragraph 'Encoding:'"
}

$newLines1 = @(
    "Quads: URN. ID Occurrence Sequence Count starts at Quad ID Count.",
    "Occurrence: Distinct URN Occurrences Sequence Count.",
    "SPOResource: First URN Occurrence ID.",
    "SPOs: (SPOResource, Occurrence);",
    "Quad : (ContextResource, Occurrence);",
    "Normalize / Aggregate CSPOs IDs (States).",
    "Graph / Tree List Parent / Child encoding / hashing:",
    "(C (S (P (O)."
)

$cursor = $anchor1
foreach ($line in $newLines1) {
    $cursor.Range.InsertParagraphAfter()
    $cursor = $cursor.Next()
    $cursor.Range.Text = $line
}
# Final empty bullet paragraph (no text run content at all).
$cursor.Range.InsertParagraphAfter()
$cursor = $cursor.Next()

# ---------------------------------------------------------------------------
# 2) Insert one new bullet paragraph right after
#    'HashedURN : "[" HashedQuad "]" | HashedCSPOString;'
# ---------------------------------------------------------------------------
$anchor2 = Find-ParagraphByText $d "HashedURN : `"[`" HashedQuad `"]`" | HashedCSPOString;`r" 1
if ($anchor2 -eq $null) {
    throw "Could not locate anchor paragraph 'HashedURN : ...'"
}

$anchor2.Range.InsertParagraphAfter()
$newPara2 = $anchor2.Next()
$newPara2.Range.Text = "HashedCSPOString : Context `":`" Subject `":`" Predicate `":`" Object;"
